$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells that look numeric stay as text (matches source formatting)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.885.64"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.73"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.00"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.89"
$ws.Range("E8").Value = "  +14.57%  "
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0692"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.115.83"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.53"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.846.22"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.75"
$ws.Range("E15").Value = "  +6.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.660"
$ws.Range("E16").Value = "  +4.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.819.48"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.38"
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.47"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.15"
$ws.Range("E21").Value = "  +8.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.74"
$ws.Range("E22").Value = "  +15.27%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.89"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.93"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.87"
$ws.Range("E27").Value = "  +3.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.123"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +7.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.97"
$ws.Range("E31").Value = "  +3.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.02"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0534"
$ws.Range("E33").Value = "  +3.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "91.23"
$ws.Range("E35").Value = "  +11.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.671"
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.343.30"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.08"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("E39").Value = "  +8.54%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.44"
$ws.Range("E40").Value = "  +2.68%  "
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.90"
$ws.Range("E42").Value = "  +8.47%  "
$ws.Range("E43").Value = "  +7.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.013.72"
$ws.Range("E47").Value = "  +2.03%  "
$ws.Range("E48").Value = "  +3.52%  "
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.11"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.27"
$ws.Range("E51").Value = "  +5.40%  "
